# Error Calculations and Plots
#
# The "missing_data.xlsx" sheet models a human-curated CDF re-impute pass:
# two rows are dropped entirely ("RM 232" and "SC 92", which shifts every
# row below them up), and a handful of individual A-F cells are revised
# (some previously-filled numbers become newly "missing", and some
# previously-missing cells get a newly imputed number).
#
# "Missing" in this sheet is represented as a literal empty-text cell
# (inlineStr with no content) rather than a truly blank cell, so we clear
# cells by assigning a bare quote-prefix ("'") — Excel's "force text, empty
# after the quote" idiom — which yields the same empty-string/Text cell
# the rest of the sheet already uses for missing values, instead of a
# genuinely blank cell. ClearFormats() afterwards drops the transient
# quote-prefix cell style so the cleared cells match the plain (unstyled)
# empty-text cells already used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the "RM 232" row (originally row 26) -----------------------------
$ws.Rows("26:26").Delete()

# --- Drop the "SC 92" row (originally row 28, now row 27 after the above) --
$ws.Rows("27:27").Delete()

# --- Per-cell value revisions on the now-shifted rows -----------------------

# RM 2: F2 18.03 -> (now missing)
$ws.Range("F2").Value = "'"

# RM 14: F5 (was missing) -> 17.66
$ws.Range("F5").Value = 17.66

# RM 21: E6 (was missing) -> -5.7 ; F6 (was missing) -> 16.43
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43

# RM 38: E8 -6.6 -> (now missing)
$ws.Range("E8").Value = "'"

# RM 42: F9 17.26 -> (now missing)
$ws.Range("F9").Value = "'"

# RM 52 a: F10 16.43 -> (now missing)
$ws.Range("F10").Value = "'"

# RM 81: E12 (was missing) -> -5.3
$ws.Range("E12").Value = -5.3

# RM 90: E14 -5.4 -> (now missing)
$ws.Range("E14").Value = "'"

# RM 116: E17 (was missing) -> -7.3
$ws.Range("E17").Value = -7.3

# RM 120: E18 (was missing) -> -8.5
$ws.Range("E18").Value = -8.5

# RM 125: E19 -6.5 -> (now missing)
$ws.Range("E19").Value = "'"

# RM 134: E20 -7.2 -> (now missing)
$ws.Range("E20").Value = "'"

# RM 140: E23 (was missing) -> -7
$ws.Range("E23").Value = -7

# RM 142a: F24 (was missing) -> 16.78
$ws.Range("F24").Value = 16.78

# SC 101 (now row 27): D27 (was missing) -> -14.6 ; E27 -10 -> (now missing)
$ws.Range("D27").Value = -14.6
$ws.Range("E27").Value = "'"

# SC 105 (now row 28): D28 -13.7 -> (now missing) ; F28 17.44 -> (now missing)
$ws.Range("D28").Value = "'"
$ws.Range("F28").Value = "'"

# SC 119 (now row 29): D29 -13 -> (now missing)
$ws.Range("D29").Value = "'"

# SC 120 (now row 30): D30 (was missing) -> -13.6 ; F30 (was missing) -> 16.89
$ws.Range("D30").Value = -13.6
$ws.Range("F30").Value = 16.89

# SC 193 (now row 32): D32 -14.7 -> (now missing)
$ws.Range("D32").Value = "'"
